$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest cryptocurrency data

$ws.Range("D2").Value = "59.243.88"
$ws.Range("E2").Value = "  +0.03%  "

$ws.Range("D3").Value = "2.513.90"
$ws.Range("E3").Value = "  -0.27%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.25"
$ws.Range("E5").Value = "  -0.93%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.28"
$ws.Range("E6").Value = "  -3.78%  "

$ws.Range("E7").Value = "  +0.33%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.564"
$ws.Range("E8").Value = "  -1.58%  "

$ws.Range("D9").Value = "2.516.51"
$ws.Range("E9").Value = "  -1.04%  "

$ws.Range("E10").Value = "  +0.36%  "

$ws.Range("E11").Value = "  +1.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.43"
$ws.Range("E12").Value = "  -3.29%  "

$ws.Range("E13").Value = "  -0.34%  "

$ws.Range("D14").Value = "2.962.89"
$ws.Range("E14").Value = "  -0.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.44"
$ws.Range("E15").Value = "  -1.18%  "

$ws.Range("D16").Value = "59.164.47"
$ws.Range("E16").Value = "  +0.01%  "

$ws.Range("E17").Value = "  +0.05%  "

$ws.Range("D18").Value = "2.520.86"
$ws.Range("E18").Value = "  -0.51%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.15"
$ws.Range("E19").Value = "  -0.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.31"
$ws.Range("E20").Value = "  +0.65%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.52"
$ws.Range("E21").Value = "  -0.04%  "

$ws.Range("E22").Value = "  +0.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.80"
$ws.Range("E23").Value = "  -0.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.65"
$ws.Range("E24").Value = "  +2.50%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.428"
$ws.Range("E25").Value = "  -0.97%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.167"
$ws.Range("E26").Value = "  +1.55%  "

$ws.Range("E27").Value = "  +0.54%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.83"
$ws.Range("E28").Value = "  -2.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.93"
$ws.Range("E29").Value = "  +3.46%  "

$ws.Range("D30").Value = "0.0₃0775"
$ws.Range("E30").Value = "  -0.76%  "

$ws.Range("E31").Value = "  -2.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "166.10"
$ws.Range("E32").Value = "  +5.35%  "

$ws.Range("E33").Value = "  +0.13%  "

$ws.Range("E34").Value = "  -2.56%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.11"
$ws.Range("E35").Value = "  -7.59%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.51"
$ws.Range("E36").Value = "  -1.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.26"
$ws.Range("E37").Value = "  -2.66%  "

$ws.Range("E38").Value = "  -1.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.86"
$ws.Range("E39").Value = "  -0.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.69"
$ws.Range("E40").Value = "  -0.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.815"
$ws.Range("E41").Value = "  -1.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.24"
$ws.Range("E42").Value = "  -6.93%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "278.97"
$ws.Range("E43").Value = "  -5.59%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  +0.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.87"
$ws.Range("E45").Value = "  +0.58%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.597"
$ws.Range("E46").Value = "  -0.68%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0930"
$ws.Range("E47").Value = "  -0.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.71"
$ws.Range("E48").Value = "  +0.43%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0513"
$ws.Range("E49").Value = "  -0.39%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0224"
$ws.Range("E50").Value = "  -1.66%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.81"
$ws.Range("E51").Value = "  -1.91%  "
